$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new week of data was recorded for "Feria Lagunitas de Puerto Montt - Poroto granado".
# Insert a new row at row 7 (this shifts the existing rows 7-29 down to 8-30, which is
# exactly what the diff shows: each old row's data reappears one row further down, and
# the dimension grows from A1:R29 to A1:R30).
$ws.Rows(7).Insert()

# Populate the newly inserted row 7 with the new observation.
$ws.Cells.Item(7, 1).Value = 4
$ws.Cells.Item(7, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(7, 3).Value = "Los Lagos"
$ws.Cells.Item(7, 4).Value = 44600
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = 100112030
$ws.Cells.Item(7, 7).Value = "Poroto granado"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 70
$ws.Cells.Item(7, 11).Value = 35000
$ws.Cells.Item(7, 12).Value = 35000
$ws.Cells.Item(7, 13).Value = 35000
$ws.Cells.Item(7, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(7, 15).Value = "Región Metropolitana"
$ws.Cells.Item(7, 16).Value = 1400
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"
